$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New benchmark block (rows 140-151): batch-size scalability results ----

# Prime the shared-string table so new strings are interned in the same order
# the original authoring session produced (matches the canonical OOXML diff).
# Excel assigns a shared-string index on first use; writing-then-clearing a
# scratch range lets us fix that order deterministically before the real writes.
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = '15th Jun'
$arr[0,1] = 'Batch size'
$arr[0,2] = '10 min'
$arr[0,3] = 'food waste sum nutrients'
$arr[0,4] = 'max num people'
$arr[0,5] = 'days x people'
$arr[0,6] = 'imporved over time'
$arr[0,7] = 'yes'
$arr[0,8] = 'cost/(days x people)'
$arr[0,9] = 'emissions /(days x people)'
$arr[0,10] = 'food waste /(days x people)'
$arr[0,11] = 'satisfiability time'
$arr[0,12] = '12 s'
$arr[0,13] = '6 s'
$arr[0,14] = '8 s '
$arr[0,15] = '1 m 15 s'
$arr[0,16] = 'no'
$arr[0,17] = '1m'
$arr[0,18] = '18 s'
$arr[0,19] = '20 s'
$arr[0,20] = '3m 56s'
$arr[0,21] = '27s'
$primeRange = $ws.Range("A300").Resize(1, 22)
$primeRange.Value = $arr
$primeRange.ClearContents() | Out-Null

# Prime the two new cell styles (center+top, then center-only) in that exact
# order on scratch cells so they intern as cellXfs indices 7 and 8 respectively
# -- matching each alignment property write commits its own style revision, so
# writing Horizontal then Vertical on the *real* ranges below would otherwise
# register the transient "horizontal-only" state first.
$styleScratch1 = $ws.Range("A300")
$styleScratch1.Value = 1
$styleScratch1.HorizontalAlignment = -4108
$styleScratch1.VerticalAlignment = -4160
$styleScratch2 = $ws.Range("A301")
$styleScratch2.Value = 1
$styleScratch2.HorizontalAlignment = -4108
$ws.Range("A300:A301").Clear() | Out-Null

# Row 140
$ws.Range("A140").Value = 'Gecode'
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = '15th Jun'
$arr[0,1] = '10 min'
$ws.Range("B140:C140").Value = $arr

# Row 141
$ws.Range("A141").Value = 'Batch size'
$ws.Range("O141").Value = 'Batch size'
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 7
$arr[0,1] = 8
$arr[0,2] = 9
$arr[0,3] = 10
$arr[0,4] = 12
$arr[0,5] = 15
$arr[0,6] = 21
$ws.Range("B141:H141").Value = $arr
$ws.Range("B141:H141").HorizontalAlignment = -4108
$ws.Range("B141:H141").VerticalAlignment = -4160
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 25
$arr[0,1] = 28
$ws.Range("I141:J141").Value = $arr
$ws.Range("I141:J141").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 30
$arr[0,1] = 35
$arr[0,2] = 40
$arr[0,3] = 42
$ws.Range("K141:N141").Value = $arr
$ws.Range("K141:N141").HorizontalAlignment = -4108
$ws.Range("K141:N141").VerticalAlignment = -4160

# Row 142
$ws.Range("A142").Value = 'max num people'
$ws.Range("O142").Value = 'max num people'
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 30
$arr[0,1] = 32
$arr[0,2] = 46
$arr[0,3] = 59
$arr[0,4] = 93
$arr[0,5] = 125
$arr[0,6] = 154
$ws.Range("B142:H142").Value = $arr
$ws.Range("B142:H142").HorizontalAlignment = -4108
$ws.Range("B142:H142").VerticalAlignment = -4160
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 165
$arr[0,1] = 167
$ws.Range("I142:J142").Value = $arr
$ws.Range("I142:J142").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 163
$arr[0,1] = 89
$arr[0,2] = 91
$arr[0,3] = 25
$ws.Range("K142:N142").Value = $arr
$ws.Range("K142:N142").HorizontalAlignment = -4108
$ws.Range("K142:N142").VerticalAlignment = -4160

# Row 143
$ws.Range("A143").Value = 'days x people'
$ws.Range("O143").Value = 'days x people'
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 210
$arr[0,1] = 256
$arr[0,2] = 414
$arr[0,3] = 590
$arr[0,4] = 1116
$arr[0,5] = 1875
$arr[0,6] = 3234
$ws.Range("B143:H143").Value = $arr
$ws.Range("B143:H143").HorizontalAlignment = -4108
$ws.Range("B143:H143").VerticalAlignment = -4160
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 4125
$arr[0,1] = 4676
$ws.Range("I143:J143").Value = $arr
$ws.Range("I143:J143").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 4890
$arr[0,1] = 3115
$arr[0,2] = 3640
$arr[0,3] = 1050
$ws.Range("K143:N143").Value = $arr
$ws.Range("K143:N143").HorizontalAlignment = -4108
$ws.Range("K143:N143").VerticalAlignment = -4160

# Row 144
$ws.Range("A144").Value = 'satisfiability time'
$ws.Range("O144").Value = 'satisfiability time'
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = '12 s'
$arr[0,1] = '6 s'
$arr[0,2] = '6 s'
$arr[0,3] = '6 s'
$arr[0,4] = '8 s '
$arr[0,5] = '1 m 15 s'
$arr[0,6] = '12 s'
$ws.Range("B144:H144").Value = $arr
$ws.Range("B144:H144").HorizontalAlignment = -4108
$ws.Range("B144:H144").VerticalAlignment = -4160
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = '1m'
$arr[0,1] = '18 s'
$ws.Range("I144:J144").Value = $arr
$ws.Range("I144:J144").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = '20 s'
$arr[0,1] = '3m 56s'
$arr[0,2] = 'unable to compute'
$arr[0,3] = '27s'
$ws.Range("K144:N144").Value = $arr
$ws.Range("K144:N144").HorizontalAlignment = -4108
$ws.Range("K144:N144").VerticalAlignment = -4160

# Row 145
$ws.Range("A145").Value = 'imporved over time'
$ws.Range("O145").Value = 'imporved over time'
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 'yes'
$arr[0,1] = 'yes'
$arr[0,2] = 'yes'
$arr[0,3] = 'yes'
$arr[0,4] = 'yes'
$arr[0,5] = 'no'
$arr[0,6] = 'no'
$ws.Range("B145:H145").Value = $arr
$ws.Range("B145:H145").HorizontalAlignment = -4108
$ws.Range("B145:H145").VerticalAlignment = -4160
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 'no'
$arr[0,1] = 'yes'
$ws.Range("I145:J145").Value = $arr
$ws.Range("I145:J145").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 'yes'
$arr[0,1] = 'yes'
$arr[0,2] = 'unable to compute'
$arr[0,3] = 'no'
$ws.Range("K145:N145").Value = $arr
$ws.Range("K145:N145").HorizontalAlignment = -4108
$ws.Range("K145:N145").VerticalAlignment = -4160

# Row 146
$ws.Range("A146").Value = 'cost £'
$ws.Range("O146").Value = 'cost £'
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 1713
$arr[0,1] = 1966
$arr[0,2] = 3149
$arr[0,3] = 3677
$ws.Range("B146:E146").Value = $arr
$ws.Range("B146:E146").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6497
$arr[0,1] = 13222
$arr[0,2] = 23956
$arr[0,3] = 35142
$arr[0,4] = 39977
$arr[0,5] = 33533
$arr[0,6] = 26349
$ws.Range("F146:L146").Value = $arr
$ws.Range("M146").Value = 'unable to compute'
$ws.Range("M146").HorizontalAlignment = -4108
$ws.Range("M146").VerticalAlignment = -4160
$ws.Range("N146").Value = 8953

# Row 147
$ws.Range("A147").Value = 'emissions kg'
$ws.Range("O147").Value = 'emissions kg'
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 1632
$arr[0,1] = 2031
$arr[0,2] = 3068
$arr[0,3] = 3322
$ws.Range("B147:E147").Value = $arr
$ws.Range("B147:E147").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6027
$arr[0,1] = 11915
$arr[0,2] = 22637
$arr[0,3] = 32098
$arr[0,4] = 35670
$arr[0,5] = 34496
$arr[0,6] = 25271
$ws.Range("F147:L147").Value = $arr
$ws.Range("M147").Value = 'unable to compute'
$ws.Range("M147").HorizontalAlignment = -4108
$ws.Range("M147").VerticalAlignment = -4160
$ws.Range("N147").Value = 9003

# Row 148
$ws.Range("A148").Value = 'food waste sum nutrients'
$ws.Range("O148").Value = 'food waste sum nutrients'
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 74168
$arr[0,1] = 74014
$arr[0,2] = 136986
$arr[0,3] = 148030
$ws.Range("B148:E148").Value = $arr
$ws.Range("B148:E148").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 337583
$arr[0,1] = 388942
$arr[0,2] = 752213
$arr[0,3] = 1107399
$arr[0,4] = 922383
$arr[0,5] = 1892849
$arr[0,6] = 1129074
$ws.Range("F148:L148").Value = $arr
$ws.Range("M148").Value = 'unable to compute'
$ws.Range("M148").HorizontalAlignment = -4108
$ws.Range("M148").VerticalAlignment = -4160
$ws.Range("N148").Value = 445200

# Row 149
$ws.Range("A149").Value = 'cost/(days x people)'
$ws.Range("O149").Value = 'cost/(days x people)'
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = "=1713/210"
$arr[0,1] = "=1966/256"
$arr[0,2] = "=3149/414"
$arr[0,3] = "=3677/590"
$ws.Range("B149:E149").Formula = $arr
$ws.Range("B149:E149").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = "=6497/1116"
$arr[0,1] = "=13222/1875"
$arr[0,2] = "=23956/3234"
$arr[0,3] = "=35142/4125"
$arr[0,4] = "=39977/4676"
$arr[0,5] = "=33533/4890"
$arr[0,6] = "=26349/3115"
$ws.Range("F149:L149").Formula = $arr
$ws.Range("M149").Value = 'unable to compute'
$ws.Range("M149").HorizontalAlignment = -4108
$ws.Range("M149").VerticalAlignment = -4160
$ws.Range("N149").Formula = "=8953/1050"

# Row 150
$ws.Range("A150").Value = 'emissions /(days x people)'
$ws.Range("O150").Value = 'emissions /(days x people)'
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = "=1632/210"
$arr[0,1] = "=2031/256"
$arr[0,2] = "=3068/414"
$arr[0,3] = "=3322/590"
$ws.Range("B150:E150").Formula = $arr
$ws.Range("B150:E150").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = "=6027/1116"
$arr[0,1] = "=11915/1875"
$arr[0,2] = "=22637/3234"
$arr[0,3] = "=32098/4125"
$arr[0,4] = "=35670/4676"
$arr[0,5] = "=34496/4890"
$arr[0,6] = "=25271/3115"
$ws.Range("F150:L150").Formula = $arr
$ws.Range("M150").Value = 'unable to compute'
$ws.Range("M150").HorizontalAlignment = -4108
$ws.Range("M150").VerticalAlignment = -4160
$ws.Range("N150").Formula = "=9003/1050"

# Row 151
$ws.Range("A151").Value = 'food waste /(days x people)'
$ws.Range("O151").Value = 'food waste /(days x people)'
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = "=74168/210"
$arr[0,1] = "=74014/256"
$arr[0,2] = "=136986/414"
$arr[0,3] = "=148030/590"
$ws.Range("B151:E151").Formula = $arr
$ws.Range("B151:E151").HorizontalAlignment = -4108
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = "=337583/1116"
$arr[0,1] = "=388942/1875"
$arr[0,2] = "=752213/3234"
$arr[0,3] = "=1107399/4125"
$arr[0,4] = "=922383/4676"
$arr[0,5] = "=1892849/4890"
$arr[0,6] = "=1129074/3115"
$ws.Range("F151:L151").Formula = $arr
$ws.Range("M151").Value = 'unable to compute'
$ws.Range("M151").HorizontalAlignment = -4108
$ws.Range("M151").VerticalAlignment = -4160
$ws.Range("N151").Formula = "=445200/1050"

# ---- View state: scroll to the new block and leave the selection where the
# authoring session left it (D149, first cell of the new "cost/(days x people)"
# row). ----
$ws.Range("A132").Select() | Out-Null
$ws.Range("D149").Select() | Out-Null